# Apply a permutation of the observation rows 4-12 (inclusive) on the
# "Artfynd" sheet. Every column (A:AY) of each row is moved to a new row
# position as captured by the mapping below (new row -> row the data
# originally came from).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 4
$lastRow = 12
$firstCol = 1   # A
$lastCol = 51   # AY

# new row -> source (old) row
$mapping = @{
    4  = 6
    5  = 11
    6  = 8
    7  = 10
    8  = 12
    9  = 5
    10 = 7
    11 = 9
    12 = 4
}

function Set-CellSafe($cell, $val) {
    if ($null -eq $val) {
        # Wholly empty cell - nothing to write.
        $cell.ClearContents()
        return
    }
    if ($val -is [string] -and $val -eq "") {
        # Preserve an explicit empty-string cell (distinct from a cell that
        # has no entry at all) the same way Excel stores a cell that was
        # typed as a lone leading apostrophe.
        $cell.Value2 = "'"
        return
    }
    $cell.Value2 = $val
    if ($val -is [string]) {
        $after = $cell.Value2
        if (-not ($after -is [string])) {
            # Excel "smart"-converted our literal string (e.g. a date- or
            # time-looking value) into a different type on assignment.
            # Force it to stay literal text using a quote prefix, same as
            # typing an apostrophe before the value in the UI.
            $cell.Value2 = "'" + $val
        }
    }
}

# 1) Snapshot every cell (Value2) for the affected rows/columns before
#    making any changes.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Write the snapshotted values back into their new row positions. Since
#    $mapping is a full permutation of 4..12, every source row is read
#    from the (already captured) snapshot, so overwriting rows in place is
#    safe regardless of iteration order.
foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $val = $snapshot["$oldRow,$c"]
        Set-CellSafe $ws.Cells.Item($newRow, $c) $val
    }
}
